$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Favorite Ice Cream" / "Favorite Pizza Toppings" values for
# row 13 (Luke Bertram), which were previously blank.
$ws.Range("E13").Value = "Blue Moon"
$ws.Range("F13").Value = "pepperoni"

# Keep the active selection consistent with the edited cell, matching
# what Excel records after the user finishes editing F13.
$ws.Range("F13").Select()
